$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data rows (anonymize attendee list) ---
# Row 2
$ws.Range("B2").Value = "A대"
$ws.Range("C2").Value = "홍길동1"
$ws.Range("D2").Value = "멘토"

# Row 3
$ws.Range("B3").Value = "C대"
$ws.Range("C3").Value = "홍길동2"
$ws.Range("D3").Value = "멘토"

# Row 4
$ws.Range("B4").Value = "D대"
$ws.Range("C4").Value = "홍길동3"
$ws.Range("D4").Value = "멘토"

# Row 5
$ws.Range("B5").Value = "E대"
$ws.Range("C5").Value = "홍길동4"
$ws.Range("D5").Value = "멘토"

# Row 6
$ws.Range("B6").Value = "B대"
$ws.Range("C6").Value = "홍길동5"
$ws.Range("D6").Value = "멘티"

# Row 7
$ws.Range("B7").Value = "D대"
$ws.Range("C7").Value = "홍길동6"
$ws.Range("D7").Value = "멘티"

# Row 8
$ws.Range("B8").Value = "B대"
$ws.Range("C8").Value = "홍길동7"
$ws.Range("D8").Value = "멘티"

# --- Mixed-font rich text for B2 ("A" in Calibri, "대" in 나눔스퀘어) ---
# (touch a scratch cell first so the Calibri font gets registered in the
# workbook's font table without altering B2's own cell style)
$ws.Range("Z1").Font.Name = "Calibri"
$ws.Range("Z1").Clear()
$ws.Range("B2").Characters(1,1).Font.Name = "Calibri"
$ws.Range("B2").Characters(2,1).Font.Name = "나눔스퀘어"

# --- Formulas now compare against "멘토" instead of "순장" ---
$ws.Range("A2").Formula = "=IF(D2=""멘토"",0,1)"
$ws.Range("A3:A8").Formula = "=IF(D3=""멘토"",0,1)"

$ws.Range("E2:E8").Formula = "=B2&"" ""&D2"
$ws.Range("E6").Formula = "=B6&"" ""&D6"
$ws.Range("E7").Formula = "=B7&"" ""&D7"

# --- Column widths (auto-fit results after content changed) ---
$ws.Columns.Item(1).ColumnWidth = 13.36
$ws.Columns.Item(3).ColumnWidth = 8.65
$ws.Columns.Item(5).ColumnWidth = 11.08

# --- Selection moved ---
$ws.Range("H10").Select()
